# Update cryptos list - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)

$ws.Range("D2").Value = "28.249.85"
$ws.Range("E2").Value = "  +5.71%  "

$ws.Range("D3").Value = "1.787.41"
$ws.Range("E3").Value = "  +3.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.52"
$ws.Range("E5").Value = "  +1.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2694"
$ws.Range("E8").Value = "  +2.66%  "

$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("D10").Value = "1.783.02"
$ws.Range("E10").Value = "  +2.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.53"
$ws.Range("E11").Value = "  +3.68%  "

$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6292"
$ws.Range("E13").Value = "  +2.92%  "

$ws.Range("E14").Value = "  +3.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "80.08"
$ws.Range("E15").Value = "  +3.70%  "

$ws.Range("D16").Value = "28.200.16"
$ws.Range("E16").Value = "  +6.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9990"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007244"
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.06"
$ws.Range("E20").Value = "  +5.64%  "

$ws.Range("D21").Value = "2.012.34"
$ws.Range("E21").Value = "  +3.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.555"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.770"
$ws.Range("E23").Value = "  +2.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.260"
$ws.Range("E24").Value = "  +3.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.05"
$ws.Range("E25").Value = "  +2.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.79"
$ws.Range("E26").Value = "  +2.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.857"
$ws.Range("E27").Value = "  +4.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "109.83"
$ws.Range("E28").Value = "  +2.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.385"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.181"
$ws.Range("E30").Value = "  +6.23%  "

$ws.Range("E31").Value = "  +3.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.769"
$ws.Range("E32").Value = "  +2.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04884"
$ws.Range("E33").Value = "  +8.83%  "

$ws.Range("E34").Value = "  +7.99%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6575"
$ws.Range("E35").Value = "  +5.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.616"
$ws.Range("E36").Value = "  +0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9474"
$ws.Range("E37").Value = "  +0.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.610"
$ws.Range("E38").Value = "  +7.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.075"
$ws.Range("E39").Value = "  +1.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.931"
$ws.Range("E40").Value = "  +6.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01550"
$ws.Range("E41").Value = "  +2.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9992"
$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.80"
$ws.Range("E43").Value = "  +0.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3997"
$ws.Range("E44").Value = "  +3.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.195"
$ws.Range("E45").Value = "  +4.00%  "

$ws.Range("E46").Value = "  +4.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05448"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.044"
$ws.Range("E48").Value = "  +2.40%  "

# Row 49/50 swap: NEARProtocol <-> Elrond
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.76"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.294"
$ws.Range("E50").Value = "  +5.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.95"
$ws.Range("E51").Value = "  +2.43%  "
